# STRAV-15 correct some column datatypes
#
# This script updates the "Staging.Activities" and "Reporting.FACT_activities"
# table definitions (and the "UNIQUEIDENTIFER" -> "UNIQUEIDENTIFIER" typo fix
# on the Landing/Staging sheets) to match the corrected column datatypes.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Landing sheet: fix "UNIQUEIDENTIFER" typo for the ID column datatype
# ---------------------------------------------------------------------
$landing = $wb.Worksheets.Item("Landing")
$landing.Range("B3").Value = "UNIQUEIDENTIFIER"

# ---------------------------------------------------------------------
# Staging sheet: fix the typo + correct several column datatypes
# ---------------------------------------------------------------------
$staging = $wb.Worksheets.Item("Staging")

# ID UNIQUEIDENTIFER -> UNIQUEIDENTIFIER
$staging.Range("B3").Value = "UNIQUEIDENTIFIER"

# athlete_id INT -> NVARCHAR(50)
$staging.Range("B4").Value = "NVARCHAR(50)"

# total_elevation_gain INT -> NUMERIC(5, 1)
$staging.Range("B9").Value = "NUMERIC(5, 1)"

# activity_id INT -> NVARCHAR(50)
$staging.Range("B12").Value = "NVARCHAR(50)"

# external_id NVARCHAR(50) NOT NULL -> NVARCHAR(100) (drop NOT NULL constraint)
$staging.Range("B13").Value = "NVARCHAR(100)"
$staging.Range("C13").Value = ""

# private BOOLEAN -> BIT NOT NULL
$staging.Range("B20").Value = "BIT NOT NULL"

# gear_id INT -> NVARCHAR(50)
$staging.Range("B21").Value = "NVARCHAR(50)"

# average_speed NUMERIC(2, 1) -> NUMERIC(3, 1)
$staging.Range("B22").Value = "NUMERIC(3, 1)"

# max_speed NUMERIC(2, 1) -> NUMERIC(3, 1)
$staging.Range("B23").Value = "NUMERIC(3, 1)"

# average_heartrate NUMERIC(3, 1) -> NUMERIC(4, 1)
$staging.Range("B24").Value = "NUMERIC(4, 1)"

# max_heartrate NUMERIC(3, 1) -> NUMERIC(4, 1)
$staging.Range("B25").Value = "NUMERIC(4, 1)"

# suffer_score INT -> NUMERIC(4, 1)
$staging.Range("B27").Value = "NUMERIC(4, 1)"

# ---------------------------------------------------------------------
# Reporting sheet: private BOOLEAN -> BIT NOT NULL
# ---------------------------------------------------------------------
$reporting = $wb.Worksheets.Item("Reporting")
$reporting.Range("B21").Value = "BIT NOT NULL"

# ---------------------------------------------------------------------
# Restore cursor / selection state on each sheet to match the author's
# final view (Staging remains the active/selected tab).
# ---------------------------------------------------------------------
$landing.Activate()
$landing.Range("C7").Select()

$reporting.Activate()
$reporting.Range("C21").Select()

$staging.Activate()
$staging.Range("C13").Select()
